# Re-apply corrected betexplorer scrape data (script run 21-11-2023 20:45).
# Several rows had their home/away + odds columns shifted to the wrong
# fixture; this restores the F:V block for each affected row, and appends
# the newly scraped Jong Almere City vs Jong Sparta Rotterdam match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 61; F = 'ADO 20 Heemskerk'; G = 2; H = 'Katwijk'; I = 1; J = 3.26; K = '07/10/2023 11:42'; L = 3.63; M = '07/10/2023 14:46'; N = 3.8; O = '07/10/2023 11:42'; P = 4.06; Q = '07/10/2023 14:46'; R = 1.88; S = '07/10/2023 11:42'; T = 1.78; U = '07/10/2023 14:46'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/ado-20-heemskerk-katwijk/SUhqgVo6/' },
    @{ Row = 62; F = 'Kozakken Boys'; G = 2; H = 'Lisse'; I = 0; J = 1.42; K = '06/10/2023 02:12'; L = 1.51; M = '07/10/2023 08:14'; N = 4.45; O = '06/10/2023 02:12'; P = 4.41; Q = '07/10/2023 13:02'; R = 5.2; S = '06/10/2023 02:12'; T = 4.96; U = '07/10/2023 12:44'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/kozakken-boys-lisse/bg9LZtoo/' },
    @{ Row = 63; F = 'Quick Boys'; G = 0; H = 'ACV Assen'; I = 3; J = 1.53; K = '07/10/2023 11:42'; L = 1.69; M = '07/10/2023 14:56'; N = 4.56; O = '07/10/2023 11:42'; P = 4.28; Q = '07/10/2023 14:57'; R = 4.34; S = '07/10/2023 11:42'; T = 3.85; U = '07/10/2023 14:57'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/quick-boys-acv-assen/jmd7cpwh/' },
    @{ Row = 76; F = 'Scheveningen'; G = 0; H = 'Katwijk'; I = 4; J = 3.32; K = '27/10/2023 02:42'; L = 2.94; M = '28/10/2023 14:26'; N = 3.55; O = '27/10/2023 02:42'; P = 3.74; Q = '28/10/2023 14:26'; R = 1.83; S = '27/10/2023 02:42'; T = 2.1; U = '28/10/2023 14:26'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/svv-scheveningen-katwijk/dSwyMKL9/' },
    @{ Row = 77; F = 'GVVV'; G = 2; H = 'Hardenberg'; I = 1; J = 2.57; K = '27/10/2023 02:42'; L = 2.53; M = '28/10/2023 14:19'; N = 3.4; O = '27/10/2023 02:42'; P = 3.58; Q = '28/10/2023 14:19'; R = 2.25; S = '27/10/2023 02:42'; T = 2.44; U = '28/10/2023 14:19'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/gvvv-hardenberg/Q3qCDdTk/' },
    @{ Row = 79; F = 'Quick Boys'; G = 3; H = 'Lisse'; I = 0; J = 1.18; K = '27/10/2023 03:12'; L = 1.19; M = '28/10/2023 14:59'; N = 6.75; O = '27/10/2023 03:12'; P = 7.41; Q = '28/10/2023 14:59'; R = 7.56; S = '27/10/2023 03:12'; T = 8.92; U = '28/10/2023 14:59'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/quick-boys-lisse/nDp8ExEq/' },
    @{ Row = 80; F = 'Kozakken Boys'; G = 1; H = 'Noordwijk'; I = 3; J = 2.26; K = '27/10/2023 03:12'; L = 2.49; M = '28/10/2023 14:58'; N = 3.34; O = '27/10/2023 03:12'; P = 3.47; Q = '28/10/2023 14:58'; R = 2.6; S = '27/10/2023 03:12'; T = 2.53; U = '28/10/2023 14:58'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/kozakken-boys-vv-noordwijk/4hsqKbjM/' },
    @{ Row = 83; F = 'Lisse'; G = 2; H = 'GVVV'; I = 2; J = 3.81; K = '03/11/2023 02:42'; L = 4.56; M = '04/11/2023 14:14'; N = 4; O = '03/11/2023 02:42'; P = 4.35; Q = '04/11/2023 14:14'; R = 1.62; S = '03/11/2023 02:42'; T = 1.57; U = '04/11/2023 14:14'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/lisse-gvvv/pSG8I6WE/' },
    @{ Row = 84; F = 'Noordwijk'; G = 0; H = 'Spakenburg'; I = 3; J = 2.42; K = '03/11/2023 02:42'; L = 2.46; M = '04/11/2023 14:22'; N = 3.44; O = '03/11/2023 02:42'; P = 3.62; Q = '04/11/2023 14:22'; R = 2.42; S = '03/11/2023 02:42'; T = 2.49; U = '04/11/2023 14:22'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/vv-noordwijk-spakenburg/6s5HGp1R/' },
    @{ Row = 85; F = 'ACV Assen'; G = 2; H = 'De Treffers'; I = 3; J = 3.41; K = '03/11/2023 02:42'; L = 3.57; M = '04/11/2023 14:09'; N = 3.78; O = '03/11/2023 02:42'; P = 3.88; Q = '04/11/2023 14:09'; R = 1.78; S = '03/11/2023 02:42'; T = 1.83; U = '04/11/2023 14:09'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/acv-assen-de-treffers/dQ2a8M87/' },
    @{ Row = 86; F = 'Excelsior Maassluis'; G = 2; H = 'Scheveningen'; I = 2; J = 2.84; K = '03/11/2023 02:42'; L = 3.42; M = '04/11/2023 14:26'; N = 3.57; O = '03/11/2023 02:42'; P = 3.65; Q = '04/11/2023 14:26'; R = 2.02; S = '03/11/2023 02:42'; T = 1.93; U = '04/11/2023 14:26'; V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/excelsior-maassluis-svv-scheveningen/IZ3e92g1/' }
)

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($u in $updates) {
    foreach ($c in $cols) {
        $ws.Range("$c$($u.Row)").Value = $u[$c]
    }
}

# --- Append the new match as row 105 ---
$newRow = 105

# Copy formatting (bold/border style on A, date format on E) from the last data row
$ws.Range("A104:V104").Copy()
$ws.Range("A$($newRow):V$($newRow)").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newValues = @{
    A = 104;
    B = 'netherlands';
    C = 'tweede-divisie';
    D = '2023-2024';
    E = 45251.8125;
    F = 'Jong Almere City';
    G = 1;
    H = 'Jong Sparta Rotterdam';
    I = 1;
    J = 2.48;
    K = '21/11/2023 07:42';
    L = 2.18;
    M = '21/11/2023 19:14';
    N = 3.67;
    O = '21/11/2023 07:42';
    P = 3.76;
    Q = '21/11/2023 19:12';
    R = 2.35;
    S = '21/11/2023 07:42';
    T = 2.77;
    U = '21/11/2023 19:14';
    V = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/jong-almere-city-jong-sparta-rotterdam/lAAEGut3/'
}
foreach ($c in $newValues.Keys) {
    $ws.Range("$c$newRow").Value = $newValues[$c]
}

